$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 4,20

# Row 2 (existing row -> updated values; D now "ECs")
$data[0,0]  = "ECs"
$data[0,1]  = "Fgf9"
$data[0,2]  = "Fgfr4"
$data[0,3]  = "ECs"
$data[0,4]  = 3
$data[0,5]  = 1
$data[0,6]  = 3.424501
$data[0,7]  = 10.273503
$data[0,8]  = 1
$data[0,9]  = 1
$data[0,10] = 2
$data[0,11] = 0.6666666666666666
$data[0,12] = 0.4279076666666666
$data[0,13] = 1.283723
$data[0,14] = 0.0198304262462706
$data[0,15] = 0.0198304262462706
$data[0,16] = 1.465370232407667
$data[0,17] = 13.188332091669
$data[0,18] = 0.0198304262462706
$data[0,19] = 0.0198304262462706

# Row 3 (new row, D = "FAPs")
$data[1,0]  = "ECs"
$data[1,1]  = "Fgf9"
$data[1,2]  = "Fgfr4"
$data[1,3]  = "FAPs"
$data[1,4]  = 3
$data[1,5]  = 1
$data[1,6]  = 3.424501
$data[1,7]  = 10.273503
$data[1,8]  = 1
$data[1,9]  = 1
$data[1,10] = 2
$data[1,11] = 0.6666666666666666
$data[1,12] = 0.296848
$data[1,13] = 0.890544
$data[1,14] = 0.01375675835913107
$data[1,15] = 0.01375675835913107
$data[1,16] = 1.016556272848
$data[1,17] = 9.149006455632
$data[1,18] = 0.01375675835913107
$data[1,19] = 0.01375675835913107

# Row 4 (new row, D = "M2")
$data[2,0]  = "ECs"
$data[2,1]  = "Fgf9"
$data[2,2]  = "Fgfr4"
$data[2,3]  = "M2"
$data[2,4]  = 3
$data[2,5]  = 1
$data[2,6]  = 3.424501
$data[2,7]  = 10.273503
$data[2,8]  = 1
$data[2,9]  = 1
$data[2,10] = 2
$data[2,11] = 0.6666666666666666
$data[2,12] = 0.07261233333333333
$data[2,13] = 0.217837
$data[2,14] = 0.003365056606611278
$data[2,15] = 0.003365056606611278
$data[2,16] = 0.2486610081123333
$data[2,17] = 2.237949073011
$data[2,18] = 0.003365056606611278
$data[2,19] = 0.003365056606611278

# Row 5 (new row, D = "sCs")
$data[3,0]  = "ECs"
$data[3,1]  = "Fgf9"
$data[3,2]  = "Fgfr4"
$data[3,3]  = "sCs"
$data[3,4]  = 3
$data[3,5]  = 1
$data[3,6]  = 3.424501
$data[3,7]  = 10.273503
$data[3,8]  = 1
$data[3,9]  = 1
$data[3,10] = 3
$data[3,11] = 1
$data[3,12] = 20.78097133333333
$data[3,13] = 62.342914
$data[3,14] = 0.9630477587879871
$data[3,15] = 0.9630477587879871
$data[3,16] = 71.16445711197133
$data[3,17] = 640.480114007742
$data[3,18] = 0.9630477587879871
$data[3,19] = 0.9630477587879871

$ws.Range("A2:T5").Value = $data
